$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4978516399860382
$ws.Range("B1").Value = 2.682286024093628
$ws.Range("C1").Value = 6.288734912872314
$ws.Range("D1").Value = 1.536264657974243
$ws.Range("E1").Value = 0.8888314366340637
